$wb = $excel.ActiveWorkbook

# Sheet1 = Overview, Sheet2 = zh-cn, Sheet3 = de-de
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Update Status text "In Translation" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Update datetime text values
$wsOverview.Range("G2").Value = "2016-08-26 10:56:19"
$wsZhCn.Range("H2").Value = "2016-08-26 10:56:15"
$wsDeDe.Range("H2").Value = "2016-08-26 10:56:19"

# Autofit the Status columns to reflect the new (longer) text width
$wsOverview.Columns.Item(5).AutoFit() | Out-Null
$wsOverview.Columns.Item(6).AutoFit() | Out-Null
$wsZhCn.Columns.Item(3).AutoFit() | Out-Null
$wsDeDe.Columns.Item(3).AutoFit() | Out-Null
